# Translations update for the "plot" survey form.
#
# Two display-string keys used by the app-designer tooling are being
# disambiguated:
#   - "settings" sheet, column "value" for the "display.title" setting
#     row is renamed to "display.title.text"
#   - "survey" sheet, the "display.text" column header (row 1) is
#     renamed to "display.prompt.text"
#
# The workbook also ends up with "survey" as the active sheet/tab, with
# cell E2 selected there, and with cell C2 selected on the "settings"
# sheet (reflecting where the editor was last working).

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"
[void]$settings.Range("C2").Select()

$survey = $wb.Worksheets.Item("survey")
$survey.Range("E1").Value = "display.prompt.text"
[void]$survey.Activate()
[void]$survey.Range("E2").Select()
